$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column E values for rows 7-14.
# Written in shared-string-creation order (E10's string is appended last)
# so the underlying sharedStrings table matches the source ordering.
$ws.Range("E7").Value  = "<arg_1> <arg_2> <arg_3> constrain"
$ws.Range("E8").Value  = "<arg_1> <arg_2> randomInt"
$ws.Range("E9").Value  = "<arg_1> <arg_2> cmpe"
$ws.Range("E11").Value = "<boolean> negate"
$ws.Range("E12").Value = "null"
$ws.Range("E13").Value = "get <variable name>"
$ws.Range("E14").Value = "set <variable name>"
$ws.Range("E10").Value = "boolean TRUE/boolean FALSE"

# Update view: scroll so column C is top-left, and selection moved to E17
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E17").Select()
